$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet / tab name (drop the "-HW20.xpc" suffix)
$ws.Name = "AlphaFiberF"

# 2. Recomputed value for G13 (Gaussian Quadrature scheme correction)
$ws.Range("G13").Value = 0.9951069915089162

# 3. Add new row 16 (index 14), mirroring row 15's formatting/label.
#    Copy A15 -> A16 and B15 -> B16 first so the number style (s="1")
#    and the shared string reference carry over correctly, then
#    overwrite the values with the new row's data.
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("B15").Copy($ws.Range("B16"))

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.068666319289838
$ws.Range("D16").Value = 1.067777675235372
$ws.Range("E16").Value = 0.9487134065268789
$ws.Range("F16").Value = 1.068666319289838
$ws.Range("G16").Value = 1.010109180398649
$ws.Range("H16").Value = 0.9379565520909054
$ws.Range("I16").Value = 0.981474051889873
$ws.Range("J16").Value = 1.067777675235372
$ws.Range("K16").Value = 1.008245540881125
$ws.Range("L16").Value = 1.038455930085482
$ws.Range("M16").Value = 1.002449530905253
